$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WST")

# Fill in previously-blank column B values (most recent period) for these rows
$ws.Range("B5").Value  = 67400000.0
$ws.Range("B7").Value  = 930800000.0
$ws.Range("B8").Value  = 211900000.0
$ws.Range("B9").Value  = 138200000.0
$ws.Range("B10").Value = 27900000.0
$ws.Range("B11").Value = 24800000.0
$ws.Range("B12").Value = 1412700000.0
$ws.Range("B14").Value = 2300000.0
$ws.Range("B16").Value = 69900000.0
$ws.Range("B18").Value = 49700000.0
$ws.Range("B19").Value = 147200000.0
$ws.Range("B22").Value = 22500000.0
$ws.Range("B23").Value = 55200000.0
$ws.Range("B24").Value = 12200000.0
$ws.Range("B26").Value = 427900000.0
$ws.Range("B28").Value = 246700000.0
$ws.Range("B29").Value = 18800000.0
$ws.Range("B30").Value = 1985400000.0
$ws.Range("B31").Value = 281700000.0
$ws.Range("B34").Value = 2739700000.0
$ws.Range("B35").Value = 73800000.0
$ws.Range("B36").Value = 1690800000.0

# Correction of an existing value
$ws.Range("C24").Value = -6000000.0

# Updated figures for Net Debt / Total Debt
$ws.Range("G37").Value = -172100000.0
$ws.Range("G38").Value = 267000000.0
